# Session 7 & test
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Marks")

# Fill in the new test scores for row 6
$ws.Range("G6").Value = 8
$ws.Range("H6").Value = 2.5

# New comment strings (order matters for shared-string table indices:
# "Not working exercise" must be added before "Passed")
$ws.Range("H7").Value = "Not working exercise"
$ws.Range("I6").Value = "Passed"

# Update the view: scroll so column D is the left-most visible column,
# and select J6 as the active cell.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("J6").Select()
